# Add a new "news" entry for the AVP part 1 blog post (id 2), right below
# the existing "hello_world" entry (row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("news")
$ws.Activate()

$ws.Range("A3").Value = 2
$ws.Range("C3").Value = "building-for-apple-vision-pro-1"
$ws.Range("D3").Value = "Short series documenting Pedro Cisdeli's hands" + [char]0x2011 + "on journey building ag" + [char]0x2011 + "tech tools for Apple Vision Pro"
$ws.Range("E3").Value = "June 27, 2025"
$ws.Range("F3").Value = "General"
$ws.Range("G3").Value = "developing_for_avp_part_1.md"
$ws.Range("H3").Value = "developing_for_avp_part_1.png"
$ws.Range("B3").Value = "Building for Apple Vision Pro: Part 1"
$ws.Range("B3").Select()
